$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D2:D51) to be stored as text so that purely
# numeric-looking values (e.g. "0.9955", "138.00") keep their exact
# original string formatting (incl. trailing zeros) instead of being
# converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '26.472.11'
$ws.Range("E2").Value = '  -0.41%  '

# Row 3
$ws.Range("D3").Value = '1.718.54'
$ws.Range("E3").Value = '  -1.41%  '

# Row 4
$ws.Range("D4").Value = '0.9955'
$ws.Range("E4").Value = '  -0.49%  '

# Row 5
$ws.Range("D5").Value = '239.82'
$ws.Range("E5").Value = '  -2.76%  '

# Row 6
$ws.Range("D6").Value = '0.9960'
$ws.Range("E6").Value = '  -0.45%  '

# Row 7
$ws.Range("D7").Value = '0.4895'
$ws.Range("E7").Value = '  -0.62%  '

# Row 8
$ws.Range("D8").Value = '0.2590'
$ws.Range("E8").Value = '  -3.24%  '

# Row 9
$ws.Range("D9").Value = '0.06180'
$ws.Range("E9").Value = '  -1.66%  '

# Row 10
$ws.Range("D10").Value = '1.723.79'
$ws.Range("E10").Value = '  -1.13%  '

# Row 11
$ws.Range("D11").Value = '0.06949'
$ws.Range("E11").Value = '  -1.30%  '

# Row 12
$ws.Range("D12").Value = '15.57'
$ws.Range("E12").Value = '  -1.01%  '

# Row 13
$ws.Range("D13").Value = '0.6030'
$ws.Range("E13").Value = '  -1.93%  '

# Row 14
$ws.Range("D14").Value = '4.457'
$ws.Range("E14").Value = '  -2.78%  '

# Row 15
$ws.Range("D15").Value = '76.57'
$ws.Range("E15").Value = '  -1.96%  '

# Row 16
$ws.Range("D16").Value = '0.9957'
$ws.Range("E16").Value = '  -0.48%  '

# Row 17
$ws.Range("D17").Value = '26.328.91'
$ws.Range("E17").Value = '  -0.99%  '

# Row 18
$ws.Range("D18").Value = '0.9955'
$ws.Range("E18").Value = '  -0.51%  '

# Row 19
$ws.Range("D19").Value = '0.000007103'
$ws.Range("E19").Value = '  -2.56%  '

# Row 20
$ws.Range("E20").Value = '  -2.36%  '

# Row 21
$ws.Range("D21").Value = '1.938.85'
$ws.Range("E21").Value = '  -1.45%  '

# Row 22
$ws.Range("D22").Value = '4.394'
$ws.Range("E22").Value = '  -4.05%  '

# Row 23
$ws.Range("D23").Value = '8.411'
$ws.Range("E23").Value = '  -3.53%  '

# Row 24
$ws.Range("D24").Value = '5.067'
$ws.Range("E24").Value = '  -3.70%  '

# Row 25
$ws.Range("D25").Value = '138.00'
$ws.Range("E25").Value = '  -1.10%  '

# Row 26
$ws.Range("E26").Value = '  -1.66%  '

# Row 27
$ws.Range("D27").Value = '1.389'
$ws.Range("E27").Value = '  -2.64%  '

# Row 28
$ws.Range("E28").Value = '  -1.32%  '

# Row 29
$ws.Range("D29").Value = '105.55'
$ws.Range("E29").Value = '  -2.00%  '

# Row 30
$ws.Range("D30").Value = '3.894'
$ws.Range("E30").Value = '  -3.67%  '

# Row 31
$ws.Range("D31").Value = '0.07915'
$ws.Range("E31").Value = '  -1.53%  '

# Row 32
$ws.Range("D32").Value = '3.615'
$ws.Range("E32").Value = '  -3.27%  '

# Row 33
$ws.Range("D33").Value = '0.04474'
$ws.Range("E33").Value = '  -3.34%  '

# Row 34
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '2.601'
$ws.Range("E34").Value = '  -0.43%  '

# Row 35
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '0.9956'
$ws.Range("E35").Value = '  -2.05%  '

# Row 36
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '0.6152'
$ws.Range("E36").Value = '  -3.66%  '

# Row 37
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").Value = '0.9483'
$ws.Range("E37").Value = '  +5.69%  '

# Row 38
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '1.994'
$ws.Range("E38").Value = '  -3.38%  '

# Row 39
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '2.384'
$ws.Range("E39").Value = '  -1.66%  '

# Row 40
$ws.Range("B40").Value = 'PaxDollar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D40").Value = '0.9954'
$ws.Range("E40").Value = '  -0.78%  '

# Row 41
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '0.01482'
$ws.Range("E41").Value = '  -1.41%  '

# Row 42
$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").Value = '99.55'
$ws.Range("E42").Value = '  -2.28%  '

# Row 43
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '5.438'
$ws.Range("E43").Value = '  +0.19%  '

# Row 44
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '0.3807'
$ws.Range("E44").Value = '  -2.81%  '

# Row 45
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").Value = '6.862'
$ws.Range("E45").Value = '  -0.13%  '

# Row 46
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").Value = '0.1155'
$ws.Range("E46").Value = '  -2.30%  '

# Row 47
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = '0.05354'
$ws.Range("E47").Value = '  -0.83%  '

# Row 48
$ws.Range("B48").Value = 'Elrond'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D48").Value = '30.44'
$ws.Range("E48").Value = '  -0.35%  '

# Row 49
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '7.741'
$ws.Range("E49").Value = '  -0.69%  '

# Row 50
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").Value = '51.14'
$ws.Range("E50").Value = '  -1.21%  '

# Row 51
$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").Value = '0.3346'
$ws.Range("E51").Value = '  -2.28%  '

# Restore the default (Normal) style on the Price column so that no
# extra cell formatting is introduced versus the original workbook.
$ws.Range("D2:D51").Style = "Normal"
